$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The "Etnia" group currently has two rows, in this order:
#   Row 16: Branco          | 303 (66%) | 511 (74%) | (empty)
#   Row 17: Pretos e Pardos | 159 (34%) | 176 (26%) | (empty)
#
# The edit reorders them (Pretos e Pardos first, then Brancos) and
# renames "Branco" -> "Brancos". Rather than physically moving table
# rows (unsupported here), we achieve the identical end state by
# swapping the text content of the two rows' cells.

$rowPretos = $t.Rows.Item(16)
$rowBrancos = $t.Rows.Item(17)

$rowPretos.Cells.Item(1).Range.Text = "Pretos e Pardos"
$rowPretos.Cells.Item(2).Range.Text = "159 (34%)"
$rowPretos.Cells.Item(3).Range.Text = "176 (26%)"

$rowBrancos.Cells.Item(1).Range.Text = "Brancos"
$rowBrancos.Cells.Item(2).Range.Text = "303 (66%)"
$rowBrancos.Cells.Item(3).Range.Text = "511 (74%)"
